$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet ("Sheet1" -> "Sheet")
$ws.Name = "Sheet"

# --- New "Description" column (F) ---
$ws.Range("F1").Value = "Description"

$desc = "29/06/2024 11:09am - 14000 paid , 6000 pending`n" + `
  "29/06/2024 09:52am - 1600 pay`n" + `
  "24/06/2024 05:03pm - register with 1000, pending 19000`n" + `
  "24/06/2024 05:03pm - register with 1000`n" + `
  "13/06/2024 02:39pm - npu`n" + `
  "01/06/2024 02:32pm - call at 6 pm today`n" + `
  "21/05/2024 06:14pm - npu`n" + `
  "09/05/2024 05:36pm - npu`n" + `
  "09/05/2024 04:32pm - call at 5:30`n" + `
  "28/04/2024 06:26pm - npu`n" + `
  "28/04/2024 06:26pm - npu`n" + `
  "16/04/2024 11:27am - npu , call again`n" + `
  "13/04/2024 10:52am - 10 yrs exp. in pp call again on 21 april`n" + `
  "13/04/2024 10:47am - reference of omkar. had a discusion , he will enroll in month end of april.`n"

$ws.Range("F2").Value = $desc

# Give F1 the same base style (bold font + thin border) that the other header
# cells already have, so it doesn't pick up a brand-new font entry.
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)   # xlPasteFormats

# Column widths (30, 10, 20, 10, 10, 80 chars). ColumnWidth vs the stored
# width has a fixed +5/6 padding offset in this engine, so subtract it to
# land exactly on the target stored widths.
$ws.Columns.Item(1).ColumnWidth = 29.166666666666668
$ws.Columns.Item(2).ColumnWidth = 9.166666666666666
$ws.Columns.Item(3).ColumnWidth = 19.166666666666668
$ws.Columns.Item(4).ColumnWidth = 9.166666666666666
$ws.Columns.Item(5).ColumnWidth = 9.166666666666666
$ws.Columns.Item(6).ColumnWidth = 79.16666666666667

# Row 2 (A2:F2) did not have a border before; give it the same thin border
# used throughout the rest of the table.
$ws.Range("A2:F2").Borders.LineStyle = 1

# --- Header row (row 1) alignment ---
# A1 & F1: bold font, thin border, wrap text only (no horizontal/vertical).
$ws.Range("A1").HorizontalAlignment = 1    # xlGeneral - clears old centering
$ws.Range("A1").VerticalAlignment = -4107  # xlBottom  - default, not serialized
$ws.Range("A1").WrapText = $true

$ws.Range("F1").HorizontalAlignment = 1
$ws.Range("F1").VerticalAlignment = -4107
$ws.Range("F1").WrapText = $true

# B1:E1: bold font, thin border, vertical-center only (no wrap).
$ws.Range("B1:E1").HorizontalAlignment = 1
$ws.Range("B1:E1").VerticalAlignment = -4108  # xlCenter
$ws.Range("B1:E1").WrapText = $false

# --- Data row (row 2) alignment ---
# A2 & F2: regular font, thin border, wrap text only.
$ws.Range("A2").HorizontalAlignment = 1
$ws.Range("A2").VerticalAlignment = -4107
$ws.Range("A2").WrapText = $true

$ws.Range("F2").HorizontalAlignment = 1
$ws.Range("F2").VerticalAlignment = -4107
$ws.Range("F2").WrapText = $true

# B2:E2: regular font, thin border, vertical-center only (no wrap).
$ws.Range("B2:E2").HorizontalAlignment = 1
$ws.Range("B2:E2").VerticalAlignment = -4108
$ws.Range("B2:E2").WrapText = $false

# Wrapping the long F2 text makes the engine auto-compute a tall custom row
# height; AutoFit the row back down so row 2 keeps its default height,
# matching the source workbook (no explicit <row ht=.../> override).
$ws.Rows.Item(2).AutoFit()
